$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear content of the D column cells (rows 3-6, 8-14) while preserving style
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()

# Update E13 and K13 with new values
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"
